$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sending-cluster label used in rows 2 and 3
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("A3").Value = "Inflammatory-Mac"

# Update row 2 values (new TPM-derived numbers)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7878926666666667
$ws.Range("H2").Value = 2.363678
$ws.Range("M2").Value = 0.05147733333333334
$ws.Range("O2").Value = 0.887188413789934
$ws.Range("P2").Value = 0.8871884137899338
$ws.Range("Q2").Value = 0.0405586134328889
$ws.Range("R2").Value = 0.365027520896
$ws.Range("S2").Value = 0.887188413789934
$ws.Range("T2").Value = 0.8871884137899338

# Update row 3 values (new TPM-derived numbers)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7878926666666667
$ws.Range("H3").Value = 2.363678
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.006545666666666668
$ws.Range("N3").Value = 0.019637
$ws.Range("O3").Value = 0.1128115862100661
$ws.Range("P3").Value = 0.1128115862100661
$ws.Range("Q3").Value = 0.005157282765111112
$ws.Range("R3").Value = 0.04641554488600001
$ws.Range("S3").Value = 0.1128115862100661
$ws.Range("T3").Value = 0.1128115862100661
